$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.254.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.843.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.62"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6740"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07434"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2950"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.86"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.855.48"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.005"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6720"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.13"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.130"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.276.09"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008326"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.45"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.51"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.193"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.08%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.75"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.709"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1403"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.79%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.508"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.175"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.070"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.194"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05301"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7613"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.873"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.136"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.674"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.330.91"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.720"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9187"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.956"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.55"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.08217"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +14.12%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.000.64"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5170"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.12"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.131"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05958"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.24%  "
